# Append 5 new daily rows (01/17/2021 .. 01/21/2021) to the "timeline" sheet,
# mirroring the existing Date/NewConfirmed/NewRecovered/NewHospitalized/
# NewDeaths/Confirmed/Recovered/Hospitalized/Deaths columns (A-I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date, NewConfirmed, NewRecovered, NewHospitalized, NewDeaths, Confirmed, Recovered, Hospitalized, Deaths
$data = @(
    @("01/17/2021", 374, 109,  264,    1, 12054, 9015, 2969, 70),
    @("01/18/2021", 369, 191,  178,    0, 12423, 9206, 3147, 70),
    @("01/19/2021", 171, 150,   21,    0, 12594, 9356, 3168, 70),
    @("01/20/2021",  59, 265, -207,    1, 12653, 9621, 2961, 71),
    @("01/21/2021", 142, 221,  -79,    0, 12795, 9842, 2882, 71)
)

$startRow = 383
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    # Leading apostrophe forces the date-like string to be stored as literal
    # text (matching the existing column A cells) instead of being parsed
    # into a date serial number. ClearFormats() then drops the quote-prefix
    # number format Excel applies, so the cell keeps the sheet's default style.
    $ws.Cells.Item($row, 1).Value = "'" + $rowData[0]
    $ws.Cells.Item($row, 1).ClearFormats()

    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
    $ws.Cells.Item($row, 8).Value = $rowData[7]
    $ws.Cells.Item($row, 9).Value = $rowData[8]
}
